$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for "Terminal Hortofrutícola Agro
# Chillán" / Berenjena. It belongs right after the header/earlier rows, at
# row 86, pushing the existing rows 86-101 down to 87-102 (the table stays
# sorted with the newest observation on top of that block).
$ws.Rows.Item(86).Insert()

# Fill in the newly inserted row 86 with the new observation.
$ws.Cells.Item(86, 1).Value  = 7
$ws.Cells.Item(86, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(86, 3).Value  = "Ñuble"
$ws.Cells.Item(86, 4).Value  = 45131
$ws.Cells.Item(86, 5).Value  = 16
$ws.Cells.Item(86, 6).Value  = 100112001
$ws.Cells.Item(86, 7).Value  = "Berenjena"
$ws.Cells.Item(86, 8).Value  = "Sin especificar"
$ws.Cells.Item(86, 9).Value  = "Primera"
$ws.Cells.Item(86, 10).Value = 50
$ws.Cells.Item(86, 11).Value = 8000
$ws.Cells.Item(86, 12).Value = 8000
$ws.Cells.Item(86, 13).Value = 8000
$ws.Cells.Item(86, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(86, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(86, 16).Value = 133
$ws.Cells.Item(86, 17).Value = 60
$ws.Cells.Item(86, 18).Value = "Hortaliza"
